$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sss"
$ws.Range("F7").Value = "ss"
$ws.Range("H8").Value = "sss"

$ws.Range("H8").Select()
